$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Ngày mượn" (D) and "Ngày trả dự kiến" (E) columns entirely,
# shifting the remaining "Tình trạng" / "Ghi chú" columns left.
$ws.Range("D1:E2").EntireColumn.Delete()
